$d = $word.ActiveDocument

# 1) "Created a web-based tool ..." achievement bullet under Belltower Books
$d.Content.Find.Execute(
    "Created a web-based tool to collect and analyze data from course book lists, class schedule, and final exam dates",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Created a tool to collect and aggregate data from course syllabi, class schedules, online prices, and final exam dates",
    2)

# 2) "Analyzed data ..." achievement bullet under Belltower Books
$d.Content.Find.Execute(
    "Analyzed data to determine the most productive buying opportunities based on book price, class size, demographics, and exam date; applied analysis to personally achieve top commission tier by volume",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Combined and analyzed data to determine the most productive buying opportunities based on book price, class size, student level, and exam date; applied analysis to personally achieve top commission tier by volume",
    2)

# 3) Job title "User Interface and Design Consultant" -> "User Experience and Design Consultant"
$d.Content.Find.Execute(
    "User Interface and Design Consultant",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "User Experience and Design Consultant",
    2)

# 4) Whitespace before the date range shrinks by 4 spaces (27 -> 23 leading spaces)
$d.Content.Find.Execute(
    "                           7/2007 - 8/2007",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "                       7/2007 - 8/2007",
    2)

# 5) "Advised the CEO ..." achievement bullet under My Big School
$d.Content.Find.Execute(
    "Advised the CEO of My Big School on Western design and user interface principles",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Advised the CEO of My Big School on design and UX principles for Western audiences",
    2)
